# Update the Lanzhou COVID-19 daily tracker with newly reported data
# (commit: "data updated on Aug. 04"). Fills in the four days that were
# previously blank placeholder rows (7/31 - 8/3/2022), matching the
# style already used by the surrounding filled rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: 2022-07-31 ---
$ws.Range("A26").Value = 44773
$ws.Range("B26").Value = 5
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0

# --- Row 27: 2022-08-01 ---
$ws.Range("A27").Value = 44774
$ws.Range("B27").Value = 5
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 0
# A27 needs the date number format that A26 already carries.
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)

# --- Row 28: 2022-08-02 ---
$ws.Range("A28").Value = 44775
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 0
# A28 needs the date number format too (rest of row 28 keeps its bold style).
$ws.Range("A26").Copy()
$ws.Range("A28").PasteSpecial(-4122)

# --- Row 29: 2022-08-03 ---
$ws.Range("A29").Value = 44776
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("A26").Copy()
$ws.Range("A29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the view/selection to reflect where the user ended up after entering
# the new rows.
[void]$ws.Range("J31").Select()

Write-Output "Filled rows 26-29 (2022-07-31 through 2022-08-03)."
